$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '69.726.11'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +4.68%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.618.99'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +4.07%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.14%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '628.84'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +4.68%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '159.15'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +6.73%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.615.55'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +3.95%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '1.00'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.33%  '
$ws.Range("E9").Value = '  +3.67%  '
$ws.Range("E10").Value = '  +7.75%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '7.37'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +7.40%  '
$ws.Range("E12").Value = '  +4.55%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000230'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +5.12%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '33.58'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +7.06%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '4.232.79'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +4.49%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '69.788.46'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +4.98%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.618.37'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +4.97%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.118'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.62%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.71'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +5.57%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '16.13'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +7.33%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '10.25'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +13.16%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '464.84'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +5.13%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.648'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +4.09%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '78.75'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +2.53%  '
$ws.Range("E25").Value = '  +9.79%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '10.78'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +6.69%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '3.764.13'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +4.59%  '
$ws.Range("E28").Value = '  +0.01%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.31'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +13.77%  '
$ws.Range("E30").Value = '  +5.76%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.73'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +10.29%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.178'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +12.30%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.61'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +7.53%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.00'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.25%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '26.59'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +4.03%  '
$ws.Range("E36").Value = '  +5.31%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.618.18'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +4.92%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '8.51'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +7.12%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.40'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +10.76%  '
$ws.Range("E40").Value = '  +0.05%  '
$ws.Range("B41").Value = 'Hedera'
$ws.Range("C41").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0927'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +7.66%  '
$ws.Range("B42").Value = 'Monero'
$ws.Range("C42").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '179.04'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +3.81%  '
$ws.Range("E43").Value = '  +0.13%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '5.70'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +4.22%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '32.29'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +20.71%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.915'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +4.00%  '
$ws.Range("E47").Value = '  +14.19%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '46.21'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +2.43%  '
$ws.Range("E49").Value = '  +10.99%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.84'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +4.07%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.269'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +9.27%  '

Write-Host "Applied all changes"